# Generate Report for Handoff
# Refresh the localization-status report with the latest handoff timestamps
# for the f94dc3ec-d597-4f3b-a9b2-c97a59e14f5c.md file across the Overview
# sheet and each per-locale (zh-cn / de-de) sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D7").Value = "2016-51-14 00:51:42"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E7").Value = "2016-03-14 00:51:38"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E7").Value = "2016-03-14 00:51:42"
